$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 0.009455353296505837
$ws.Range("E2").Value = 0.009455353296505837

$ws.Range("D3").Value = 0.0009150812438297632
$ws.Range("E3").Value = 0.0009150812438297632

$ws.Range("D4").Value = 1.688933277629765 / 100000000
$ws.Range("E4").Value = 1.688933277629765 / 100000000

$ws.Range("D5").Value = 6.016853034582032 / 100000000000
$ws.Range("E5").Value = 6.016853034582032 / 100000000000

$ws.Range("D6").Value = 0.001715415848908873
$ws.Range("E6").Value = 0.001715415848908873

$ws.Range("D7").Value = 0.9999999946319897
$ws.Range("E7").Value = 5.368010258521849 / 1000000000

$ws.Range("D8").Value = 0.02650235480815785
$ws.Range("E8").Value = 0.9734976451918421

$ws.Range("D9").Value = 0.7004916149084555
$ws.Range("E9").Value = 0.2995083850915445

$ws.Range("D10").Value = 0.9999999999999998
$ws.Range("E10").Value = 2.220446049250313 / 10000000000000000

$ws.Range("D11").Value = 0.9999998762598301
$ws.Range("E11").Value = 1.237401698883644 / 10000000
$ws.Range("F11").Value = 0.3998626172542572
